$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.571.85'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.408.08'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.40%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '507.90'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -2.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.37'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +4.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.558'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.449.09'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.18%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0980'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.323'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.62'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -6.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.885.59'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '57.867.06'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.87'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000135'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.09%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.445.24'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.31'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.33%  '
$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '315.61'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.10'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.43'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +5.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.41'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.07%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.155'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -4.64%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.57'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +5.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.07'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0734'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.12'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.19%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.992'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '18.07'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.23'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +4.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.88'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +3.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.79'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.809'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.94%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '135.47'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +13.46%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.39'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +1.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.91'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +3.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.573'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '255.49'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0917'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0493'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0214'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.25%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.22'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +1.91%  '
